$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the existing "Total" row from row 17 down to row 19, ---
# --- to make room for two new component rows (16 and 17).      ---

# Copy the Total row's formatting (bold font + number formats) to its new home
$ws.Range("A17:E17").Copy() | Out-Null
$ws.Range("A19:E19").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Re-enter the Total row's label and SUM formulas at row 19
$ws.Range("A19").Value = "Total"
$ws.Range("C19").Formula = "=SUM(C2:C15)"
$ws.Range("D19").Formula = "=SUM(D2:D10,D12:D14,13)"
$ws.Range("E19").Formula = "=SUM(E2:E15)"

# Clear the old Total row (row 17) completely - both contents and formatting -
# so the row becomes available again for new part data
$ws.Range("A17:E17").Clear() | Out-Null

# --- Add the two new misc. electronic components ---
# (the diode's name is entered first so it lands at the lower shared-string index,
#  matching the order the parts were added in the source workbook)
$ws.Range("A17").Value = "1N4004 diode"
$ws.Range("D17").Value = 1

$ws.Range("A16").Value = "200 Ω resistor"
$ws.Range("D16").Value = 1

# --- Restore the active cell selection ---
$ws.Range("D24").Select() | Out-Null
